# fix: changed the allocations
$wb = $excel.ActiveWorkbook

# The "allocation" sheet holds the values being corrected.
$ws = $wb.Worksheets.Item("allocation")

# Update the allocation percentages.
$ws.Range("C2").Value = 0.35
$ws.Range("C3").Value = 0.3
$ws.Range("C4").Value = 0.35

# Move the active selection on that sheet (matches the recorded cursor move).
$ws.Activate()
$ws.Range("H8").Select()
